$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Serpinf1"
$ws.Range("C2").Value = "Plxdc2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 2.063913333333333
$ws.Range("H2").Value = 6.19174
$ws.Range("I2").Value = 0.003360194407478493
$ws.Range("J2").Value = 0.003360194407478493
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 1.081988
$ws.Range("N2").Value = 3.245964
$ws.Range("O2").Value = 0.008303622335279906
$ws.Range("P2").Value = 0.008303622335279906
$ws.Range("Q2").Value = 2.233129459706666
$ws.Range("R2").Value = 20.09816513736
$ws.Range("S2").Value = 0.00002790178533282104
$ws.Range("T2").Value = 0.00002790178533282104
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Serpinf1"
$ws.Range("C3").Value = "Plxdc2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 2.063913333333333
$ws.Range("H3").Value = 6.19174
$ws.Range("I3").Value = 0.003360194407478493
$ws.Range("J3").Value = 0.003360194407478493
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 113.3348083333333
$ws.Range("N3").Value = 340.004425
$ws.Range("O3").Value = 0.8697780805714425
$ws.Range("P3").Value = 0.8697780805714423
$ws.Range("Q3").Value = 233.9132220499444
$ws.Range("R3").Value = 2105.2189984495
$ws.Range("S3").Value = 0.002922623442083539
$ws.Range("T3").Value = 0.002922623442083538
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Serpinf1"
$ws.Range("C4").Value = "Plxdc2"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 2.063913333333333
$ws.Range("H4").Value = 6.19174
$ws.Range("I4").Value = 0.003360194407478493
$ws.Range("J4").Value = 0.003360194407478493
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 15.886336
$ws.Range("N4").Value = 47.659008
$ws.Range("O4").Value = 0.1219182970932776
$ws.Range("P4").Value = 0.1219182970932776
$ws.Range("Q4").Value = 32.78802068821333
$ws.Range("R4").Value = 295.09218619392
$ws.Range("S4").Value = 0.0004096691800621328
$ws.Range("T4").Value = 0.0004096691800621328
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Serpinf1"
$ws.Range("C5").Value = "Plxdc2"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 585.4210406666667
$ws.Range("H5").Value = 1756.263122
$ws.Range("I5").Value = 0.9531061576560254
$ws.Range("J5").Value = 0.9531061576560252
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 1.081988
$ws.Range("N5").Value = 3.245964
$ws.Range("O5").Value = 0.008303622335279906
$ws.Range("P5").Value = 0.008303622335279906
$ws.Range("Q5").Value = 633.4185409488454
$ws.Range("R5").Value = 5700.766868539608
$ws.Range("S5").Value = 0.007914233578605384
$ws.Range("T5").Value = 0.007914233578605384
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Serpinf1"
$ws.Range("C6").Value = "Plxdc2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 585.4210406666667
$ws.Range("H6").Value = 1756.263122
$ws.Range("I6").Value = 0.9531061576560254
$ws.Range("J6").Value = 0.9531061576560252
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 113.3348083333333
$ws.Range("N6").Value = 340.004425
$ws.Range("O6").Value = 0.8697780805714425
$ws.Range("P6").Value = 0.8697780805714423
$ws.Range("Q6").Value = 66348.58143825721
$ws.Range("R6").Value = 597137.2329443148
$ws.Range("S6").Value = 0.8289908443868803
$ws.Range("T6").Value = 0.8289908443868802
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Serpinf1"
$ws.Range("C7").Value = "Plxdc2"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 585.4210406666667
$ws.Range("H7").Value = 1756.263122
$ws.Range("I7").Value = 0.9531061576560254
$ws.Range("J7").Value = 0.9531061576560252
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 15.886336
$ws.Range("N7").Value = 47.659008
$ws.Range("O7").Value = 0.1219182970932776
$ws.Range("P7").Value = 0.1219182970932776
$ws.Range("Q7").Value = 9300.195353500332
$ws.Range("R7").Value = 83701.75818150298
$ws.Range("S7").Value = 0.1162010796905396
$ws.Range("T7").Value = 0.1162010796905396
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Serpinf1"
$ws.Range("C8").Value = "Plxdc2"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 26.739428
$ws.Range("H8").Value = 80.218284
$ws.Range("I8").Value = 0.04353364793649628
$ws.Range("J8").Value = 0.04353364793649627
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 1.081988
$ws.Range("N8").Value = 3.245964
$ws.Range("O8").Value = 0.008303622335279906
$ws.Range("P8").Value = 0.008303622335279906
$ws.Range("Q8").Value = 28.931740222864
$ws.Range("R8").Value = 260.385662005776
$ws.Range("S8").Value = 0.0003614869713417025
$ws.Range("T8").Value = 0.0003614869713417024
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Serpinf1"
$ws.Range("C9").Value = "Plxdc2"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 26.739428
$ws.Range("H9").Value = 80.218284
$ws.Range("I9").Value = 0.04353364793649628
$ws.Range("J9").Value = 0.04353364793649627
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 113.3348083333333
$ws.Range("N9").Value = 340.004425
$ws.Range("O9").Value = 0.8697780805714425
$ws.Range("P9").Value = 0.8697780805714423
$ws.Range("Q9").Value = 3030.507947322967
$ws.Range("R9").Value = 27274.5715259067
$ws.Range("S9").Value = 0.03786461274247867
$ws.Range("T9").Value = 0.03786461274247866
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Serpinf1"
$ws.Range("C10").Value = "Plxdc2"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 26.739428
$ws.Range("H10").Value = 80.218284
$ws.Range("I10").Value = 0.04353364793649628
$ws.Range("J10").Value = 0.04353364793649627
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 15.886336
$ws.Range("N10").Value = 47.659008
$ws.Range("O10").Value = 0.1219182970932776
$ws.Range("P10").Value = 0.1219182970932776
$ws.Range("Q10").Value = 424.791537655808
$ws.Range("R10").Value = 3823.123838902272
$ws.Range("S10").Value = 0.005307548222675905
$ws.Range("T10").Value = 0.005307548222675904
